$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$preprocess = 'space after punctuation, convert unicode to ascii, trim "space" and ",", remove break line, convert to lower, remove multiple spaces'
$features = '14 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), #digit/#ascii, %kwName, %kwAddress, %kwPhone, #max_digit_skip_0 >= 7, #max_digit_skip_0 = 0, #max_ascii_skip_0 >= 7, #max_ascii_skip_0 = 0, first_character_ascii, first_character_digit, last_character_ascii, last_character_digit'
$model = 'Neuron Network'
$modelDetails = '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000'
$templateFilter = '0 filters: '

$rows = @(
    @{ r = 13; Time = '20160415_171720'; RunningTime = 1013.377; Test = 0.992666666666667; Val = 0.947194719471947; J = 0.0238095238095238 },
    @{ r = 14; Time = '20160415_173413'; RunningTime = 959.961;  Test = 0.991333333333333; Val = 0.940594059405941; J = 0.0365853658536585 },
    @{ r = 15; Time = '20160415_175013'; RunningTime = 956.226;  Test = 0.99;               Val = 0.940594059405941; J = 0.024390243902439 },
    @{ r = 16; Time = '20160415_180610'; RunningTime = 955.222;  Test = 0.988;               Val = 0.943894389438944; J = 0.036144578313253 },
    @{ r = 17; Time = '20160415_182205'; RunningTime = 970.493;  Test = 0.995333333333333; Val = 0.940594059405941; J = 0.024390243902439 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.Time
    $ws.Cells.Item($r, 2).Value = $row.RunningTime
    $ws.Cells.Item($r, 3).Value = $preprocess
    $ws.Cells.Item($r, 4).Value = $features
    $ws.Cells.Item($r, 5).Value = $model
    $ws.Cells.Item($r, 6).Value = $modelDetails
    $ws.Cells.Item($r, 7).Value = $row.Test
    $ws.Cells.Item($r, 8).Value = $row.Val
    $ws.Cells.Item($r, 9).Value = $templateFilter
    $ws.Cells.Item($r, 10).Value = $row.J
}
